# Updated cryptos list values (Price / Volume(1h) / swapped rows) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E in this sheet are stored as text; force text format so numeric-looking
# strings like "595.12" are not auto-converted to numbers by COM value assignment.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.213.71"
$ws.Range("E2").Value = "  -4.74%  "

$ws.Range("D3").Value = "3.243.34"
$ws.Range("E3").Value = "  -7.88%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "595.12"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "153.93"
$ws.Range("E6").Value = "  -11.86%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "3.235.15"
$ws.Range("E8").Value = "  -8.00%  "

$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  -11.27%  "

$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  -10.44%  "

$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  -11.36%  "

$ws.Range("D12").Value = "0.490"
$ws.Range("E12").Value = "  -16.00%  "

$ws.Range("D13").Value = "39.14"
$ws.Range("E13").Value = "  -15.32%  "

$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -11.87%  "

$ws.Range("D15").Value = "3.754.47"
$ws.Range("E15").Value = "  -8.21%  "

$ws.Range("D16").Value = "67.141.21"
$ws.Range("E16").Value = "  -4.88%  "

$ws.Range("D17").Value = "3.226.18"
$ws.Range("E17").Value = "  -8.22%  "

$ws.Range("E18").Value = "  -4.59%  "

$ws.Range("D19").Value = "532.60"
$ws.Range("E19").Value = "  -12.44%  "

$ws.Range("D20").Value = "7.04"
$ws.Range("E20").Value = "  -15.05%  "

$ws.Range("D21").Value = "14.79"
$ws.Range("E21").Value = "  -15.26%  "

$ws.Range("D22").Value = "0.757"
$ws.Range("E22").Value = "  -13.75%  "

$ws.Range("D23").Value = "7.94"
$ws.Range("E23").Value = "  -12.26%  "

$ws.Range("D24").Value = "85.63"
$ws.Range("E24").Value = "  -13.01%  "

$ws.Range("D25").Value = "13.56"
$ws.Range("E25").Value = "  -13.09%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "3.21"
$ws.Range("E27").Value = "  -13.84%  "

$ws.Range("D28").Value = "2.17"
$ws.Range("E28").Value = "  -15.55%  "

$ws.Range("D29").Value = "8.10"
$ws.Range("E29").Value = "  -10.37%  "

$ws.Range("D30").Value = "29.19"
$ws.Range("E30").Value = "  -13.84%  "

$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  -10.97%  "

$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -9.04%  "

$ws.Range("D33").Value = "535.48"
$ws.Range("E33").Value = "  -15.93%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  -15.82%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "6.50"
$ws.Range("E35").Value = "  -19.27%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").Value = "53.81"
$ws.Range("E37").Value = "  -5.30%  "

$ws.Range("D38").Value = "0.0428"
$ws.Range("E38").Value = "  -9.68%  "

$ws.Range("D39").Value = "0.0862"
$ws.Range("E39").Value = "  -13.30%  "

$ws.Range("D40").Value = "9.28"
$ws.Range("E40").Value = "  -13.80%  "

$ws.Range("D41").Value = "0.124"
$ws.Range("E41").Value = "  -13.15%  "

$ws.Range("D42").Value = "2.77"
$ws.Range("E42").Value = "  -22.80%  "

$ws.Range("D43").Value = "2.929.61"
$ws.Range("E43").Value = "  -13.15%  "

$ws.Range("D44").Value = "0.266"
$ws.Range("E44").Value = "  -14.30%  "

$ws.Range("D45").Value = "0.0₃0593"
$ws.Range("E45").Value = "  -20.08%  "

$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -15.59%  "

$ws.Range("D47").Value = "2.15"
$ws.Range("E47").Value = "  -15.76%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "26.29"
$ws.Range("E48").Value = "  -18.20%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").Value = "0.115"
$ws.Range("E50").Value = "  -11.64%  "

$ws.Range("D51").Value = "119.36"
$ws.Range("E51").Value = "  -10.59%  "
